$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Level" column (column B) entirely - data for IconX/IconY/IconWidth/
# IconHeight/Icon shifts left by one column, and the now-unused "Level"/"\u7b49\u7ea7"
# shared strings are dropped.
$ws.Range("B1").EntireColumn.Delete()

# The structural column removal above does not resync the worksheet Table
# definition (interior ListColumn deletion is not supported directly), so
# rebuild the table over the new A1:F20 range with the correct column names
# and restore its original style.
$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$tableStyle = $lo.TableStyle.Name
$lo.Unlist()
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:F20"), $null, 1)
$lo2.Name = $tableName
$lo2.TableStyle = $tableStyle

# Move the active selection from the old M9 to E9 (within the new, narrower
# table bounds).
$ws.Range("E9").Select()
